# Repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows whose data changed on re-pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 8
    4  = -1
    8  = 3
    9  = -2
    12 = -7
    17 = -1
    20 = 3
    23 = 8
    25 = 2
    30 = -7
    33 = -1
    34 = 0
    36 = -1
    41 = -1
    42 = -2
    43 = -2
    44 = 2
    46 = -1
    48 = 4
    49 = 4
    53 = -2
    54 = 4
    56 = -4
    57 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
